$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.NumberFormat = "General"
    $r.Style = "Normal"
}

Set-TextValue 'D2' '65.383.35'
Set-TextValue 'E2' '  +2.50%  '
Set-TextValue 'D3' '3.202.29'
Set-TextValue 'E3' '  +1.86%  '
Set-TextValue 'E4' '  +0.09%  '
Set-TextValue 'D5' '597.40'
Set-TextValue 'D6' '156.32'
Set-TextValue 'E6' '  +7.07%  '
Set-TextValue 'E7' '  +0.09%  '
Set-TextValue 'D8' '3.202.19'
Set-TextValue 'E8' '  +1.74%  '
Set-TextValue 'E9' '  +2.36%  '
Set-TextValue 'D10' '0.167'
Set-TextValue 'E10' '  +3.75%  '
Set-TextValue 'E11' '  +5.83%  '
Set-TextValue 'D12' '0.473'
Set-TextValue 'E12' '  +3.15%  '
Set-TextValue 'D13' '0.0000256'
Set-TextValue 'E13' '  +3.44%  '
Set-TextValue 'D14' '39.61'
Set-TextValue 'E14' '  +6.68%  '
Set-TextValue 'D15' '3.731.68'
Set-TextValue 'E15' '  +1.82%  '
Set-TextValue 'E16' '  +0.44%  '
Set-TextValue 'D17' '7.46'
Set-TextValue 'E17' '  +5.21%  '
Set-TextValue 'D18' '65.063.20'
Set-TextValue 'D19' '3.201.80'
Set-TextValue 'E19' '  +1.70%  '
Set-TextValue 'D20' '484.46'
Set-TextValue 'E20' '  +4.33%  '
Set-TextValue 'D21' '15.13'
Set-TextValue 'E21' '  +5.80%  '
Set-TextValue 'D22' '0.772'
Set-TextValue 'E22' '  +5.61%  '
Set-TextValue 'D23' '7.92'
Set-TextValue 'E23' '  +6.36%  '
Set-TextValue 'D24' '13.81'
Set-TextValue 'E24' '  +6.33%  '
Set-TextValue 'D25' '2.46'
Set-TextValue 'E25' '  +11.55%  '
Set-TextValue 'D26' '83.70'
Set-TextValue 'E26' '  +2.90%  '
Set-TextValue 'B27' 'RenderToken'
Set-TextValue 'C27' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D27' '10.09'
Set-TextValue 'E27' '  +9.84%  '
Set-TextValue 'B28' 'Dai'
Set-TextValue 'C28' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 'D28' '1.00'
Set-TextValue 'E28' '  +0.39%  '
Set-TextValue 'D29' '2.79'
Set-TextValue 'E29' '  +4.11%  '
Set-TextValue 'B30' 'NEARProtocol'
Set-TextValue 'C30' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 'D30' '7.55'
Set-TextValue 'E30' '  +8.26%  '
Set-TextValue 'B31' 'ImmutableX'
Set-TextValue 'C31' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D31' '2.29'
Set-TextValue 'E31' '  +3.24%  '
Set-TextValue 'E32' '  +0.21%  '
Set-TextValue 'E33' '  +9.09%  '
Set-TextValue 'D34' '28.56'
Set-TextValue 'E34' '  +5.70%  '
Set-TextValue 'D35' '0.0₃0906'
Set-TextValue 'E35' '  +5.82%  '
Set-TextValue 'D36' '3.62'
Set-TextValue 'E36' '  +7.90%  '
Set-TextValue 'D37' '1.10'
Set-TextValue 'E37' '  +5.04%  '
Set-TextValue 'B38' 'Stacks'
Set-TextValue 'C38' 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D38' '2.40'
Set-TextValue 'E38' '  +4.32%  '
Set-TextValue 'B39' 'Filecoin'
Set-TextValue 'C39' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D39' '6.36'
Set-TextValue 'E39' '  +5.73%  '
Set-TextValue 'D40' '478.94'
Set-TextValue 'E40' '  +8.49%  '
Set-TextValue 'D41' '9.50'
Set-TextValue 'E41' '  +8.04%  '
Set-TextValue 'D42' '51.86'
Set-TextValue 'E42' '  +1.66%  '
Set-TextValue 'D43' '0.304'
Set-TextValue 'E43' '  +9.26%  '
Set-TextValue 'E44' '  +3.35%  '
Set-TextValue 'D45' '2.964.69'
Set-TextValue 'E45' '  +1.80%  '
Set-TextValue 'E46' '  +4.27%  '
Set-TextValue 'D47' '38.97'
Set-TextValue 'E47' '  +5.47%  '
Set-TextValue 'D48' '132.39'
Set-TextValue 'E48' '  +5.28%  '
Set-TextValue 'D49' '2.36'
Set-TextValue 'E49' '  +8.52%  '
Set-TextValue 'D50' '26.00'
Set-TextValue 'E50' '  +6.71%  '
